# Refresh the cryptocurrency price/volume snapshot in columns D (Price)
# and E (Volume(1h)) for rows 2-51, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.853.80"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.598.20"
$ws.Range("E3").Value = "  +3.17%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.20"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.84"
$ws.Range("E6").Value = "  -3.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.590.78"
$ws.Range("E7").Value = "  +3.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.194"
$ws.Range("E10").Value = "  -6.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.61"
$ws.Range("E11").Value = "  +20.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.597"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.80"
$ws.Range("E13").Value = "  -3.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000279"
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.218.28"
$ws.Range("E15").Value = "  +4.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "659.54"
$ws.Range("E16").Value = "  -4.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.74"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.638.94"
$ws.Range("E18").Value = "  +4.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.215.75"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.121"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.51"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.22"
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.918"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.87"
$ws.Range("E24").Value = "  -2.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.78"
$ws.Range("E25").Value = "  -3.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  -2.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.74"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.85"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.11"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.29"
$ws.Range("E31").Value = "  -4.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.83"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.37"
$ws.Range("E33").Value = "  -6.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.21"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.89"
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "570.34"
$ws.Range("E36").Value = "  -1.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.92"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.105"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.69"
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.551.16"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0447"
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.139"
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.338"
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.23"
$ws.Range("E45").Value = "  -4.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0730"
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.65"
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"
$ws.Range("E48").Value = "  +5.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.131"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.26"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.89"
$ws.Range("E51").Value = "  +4.59%  "
